$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the rows that are no longer part of the report.
# (Delete from the bottom up so earlier row numbers stay valid.)
$ws.Rows(21).Delete()   # "Test Lead/Manager"
$ws.Rows(20).Delete()   # "- Comments/Notes"
$ws.Rows(2).Delete()    # "Test Phase"

# After the deletions the remaining rows are numbered 1-19; fill in the
# values for the test summary report.
$ws.Range("B1").Value = "Stranger Team PayAPI"
$ws.Range("B2").Value = 45257
$ws.Range("B3").Value = 45268
$ws.Range("B4").Value = 20
$ws.Range("B5").Value = 7
$ws.Range("B6").Value = 13

$ws.Range("B8").Value = "Windows 11 PRO 23H2 22631.2715"

$browsers = "Mozilla Firefox 119.0.1 (64-bit)" + [char]10 + `
  "Firefox 119.0(64-bit)" + [char]9 + [char]10 + `
  "Google Chrome Version 119.0.6045.160 (Official Build) (64-bit)" + [char]9 + [char]9 + [char]10 + `
  "Microsoft Edge Version 119.0.2151.72 (Official build) (64-bit)" + [char]9 + [char]9 + [char]10 + `
  "Opera One(version: 105.0.4970.21)" + [char]9 + [char]10 + `
  "Mobile" + [char]9
$ws.Range("B9").Value = $browsers
$ws.Rows(9).RowHeight = 180

$ws.Range("B10").Value = "-"

$ws.Range("B12").Value = 13

$ws.Range("B14").Value = "-"
$ws.Range("B15").Value = 3
$ws.Range("B16").Value = 10

$ws.Range("B18").Value = "Pass"

$ws.Range("B19").Value = 45271

$ws.Range("G14").Select()
